$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027602551808685
$ws.Range("D2").Value = 1.037116872136336
$ws.Range("E2").Value = 1.027676401288683
$ws.Range("F2").Value = 1.04803363297722
$ws.Range("I2").Value = 1.035668660751597
$ws.Range("J2").Value = 1.032759511240333
$ws.Range("K2").Value = 1.03990867375309
$ws.Range("L2").Value = 1.030495435620602
$ws.Range("M2").Value = 1.050794606922874
$ws.Range("N2").Value = 1.034226148205466

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028551203023395
$ws.Range("D3").Value = 1.037839589069555
$ws.Range("E3").Value = 1.028481520816063
$ws.Range("F3").Value = 1.048916043028126
$ws.Range("I3").Value = 1.035859814532265
$ws.Range("J3").Value = 1.03334842150747
$ws.Range("K3").Value = 1.040441362746457
$ws.Range("L3").Value = 1.031108337623683
$ws.Range("M3").Value = 1.051488798537796
$ws.Range("N3").Value = 1.034815894792731

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029165399684953
$ws.Range("D4").Value = 1.038307195305653
$ws.Range("E4").Value = 1.029003185681467
$ws.Range("F4").Value = 1.049487262765823
$ws.Range("I4").Value = 1.03598183107993
$ws.Range("J4").Value = 1.033729242753045
$ws.Range("K4").Value = 1.040785339311255
$ws.Range("L4").Value = 1.031504974339838
$ws.Range("M4").Value = 1.051937591323802
$ws.Range("N4").Value = 1.035197256848155

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029423692226647
$ws.Range("D5").Value = 1.038503765963702
$ws.Range("E5").Value = 1.02922265951806
$ws.Range("F5").Value = 1.049727459834295
$ws.Range("I5").Value = 1.036032725668431
$ws.Range("J5").Value = 1.033889280833156
$ws.Range("K5").Value = 1.040929776257755
$ws.Range("L5").Value = 1.031671730900404
$ws.Range("M5").Value = 1.052126167782445
$ws.Range("N5").Value = 1.035357522200698

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029467065595119
$ws.Range("D6").Value = 1.03853677035491
$ws.Range("E6").Value = 1.029259519831951
$ws.Range("F6").Value = 1.049767793236004
$ws.Range("I6").Value = 1.036041247548307
$ws.Range("J6").Value = 1.033916148432676
$ws.Range("K6").Value = 1.040954017808551
$ws.Range("L6").Value = 1.031699730645295
$ws.Range("M6").Value = 1.052157824930029
$ws.Range("N6").Value = 1.035384427955291

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029168850672603
$ws.Range("D7").Value = 1.038309821937463
$ws.Range("E7").Value = 1.029006117651323
$ws.Range("F7").Value = 1.049490472071527
$ws.Range("I7").Value = 1.03598251271313
$ws.Range("J7").Value = 1.033731381423441
$ws.Range("K7").Value = 1.040787269957333
$ws.Range("L7").Value = 1.031507202509272
$ws.Range("M7").Value = 1.051940111470665
$ws.Range("N7").Value = 1.035199398555708

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027923078675807
$ws.Range("D8").Value = 1.037361125114477
$ws.Range("E8").Value = 1.027948349727456
$ws.Range("F8").Value = 1.048331796690489
$ws.Range("I8").Value = 1.035733607988922
$ws.Range("J8").Value = 1.032958585965439
$ws.Range("K8").Value = 1.040088844598032
$ws.Range("L8").Value = 1.030702558131276
$ws.Range("M8").Value = 1.051029293350873
$ws.Range("N8").Value = 1.034425505639516

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025730636079627
$ws.Range("D9").Value = 1.03568915746836
$ws.Range("E9").Value = 1.026089837731969
$ws.Range("F9").Value = 1.046291977983634
$ws.Range("I9").Value = 1.035282226221023
$ws.Range("J9").Value = 1.031595002958767
$ws.Range("K9").Value = 1.038852749587892
$ws.Range("L9").Value = 1.029285086295785
$ws.Range("M9").Value = 1.049421343521029
$ws.Range("N9").Value = 1.033059986188569

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024270921737622
$ws.Range("D10").Value = 1.03457443359874
$ws.Range("E10").Value = 1.024854546031448
$ws.Range("F10").Value = 1.044933474759635
$ws.Range("I10").Value = 1.034972754085468
$ws.Range("J10").Value = 1.030684781854986
$ws.Range("K10").Value = 1.038025133298
$ws.Range("L10").Value = 1.028340444417391
$ws.Range("M10").Value = 1.048347455813991
$ws.Range("N10").Value = 1.032148472466417

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023639314108983
$ws.Range("D11").Value = 1.034091744653854
$ws.Range("E11").Value = 1.024320548882736
$ws.Range("F11").Value = 1.044345572414097
$ws.Range("I11").Value = 1.034836729579663
$ws.Range("J11").Value = 1.030290381273719
$ws.Range("K11").Value = 1.037665936875202
$ws.Range("L11").Value = 1.027931497082066
$ws.Range("M11").Value = 1.047882010788797
$ws.Range("N11").Value = 1.031753511791086

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023404776576478
$ws.Range("D12").Value = 1.033912452937113
$ws.Range("E12").Value = 1.024122333916378
$ws.Range("F12").Value = 1.044127251588689
$ws.Range("I12").Value = 1.034785901027688
$ws.Range("J12").Value = 1.030143843851116
$ws.Range("K12").Value = 1.03753239119247
$ws.Range("L12").Value = 1.027779610157541
$ws.Range("M12").Value = 1.047709058226258
$ws.Range("N12").Value = 1.031606766268535

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023455082507987
$ws.Range("D13").Value = 1.03395091157665
$ws.Range("E13").Value = 1.024164845553898
$ws.Range("F13").Value = 1.044174079758151
$ws.Range("I13").Value = 1.034796817625453
$ws.Range("J13").Value = 1.030175278398034
$ws.Range("K13").Value = 1.037561042795806
$ws.Range("L13").Value = 1.027812189754559
$ws.Range("M13").Value = 1.047746160084809
$ws.Range("N13").Value = 1.031638245456115

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023619925737604
$ws.Range("D14").Value = 1.034076924327531
$ws.Range("E14").Value = 1.024304161603207
$ws.Range("F14").Value = 1.044327524871803
$ws.Range("I14").Value = 1.034832534252978
$ws.Range("J14").Value = 1.030278269244937
$ws.Range("K14").Value = 1.03765490047386
$ws.Range("L14").Value = 1.027918941763645
$ws.Range("M14").Value = 1.047867715806228
$ws.Range("N14").Value = 1.031741382561834

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023721500347082
$ws.Range("D15").Value = 1.034154565042806
$ws.Range("E15").Value = 1.024390016793088
$ws.Range("F15").Value = 1.04442207444106
$ws.Range("I15").Value = 1.03485450031553
$ws.Range("J15").Value = 1.030341720099557
$ws.Range("K15").Value = 1.037712712868161
$ws.Range("L15").Value = 1.027984717127024
$ws.Range("M15").Value = 1.047942601648049
$ws.Range("N15").Value = 1.031804923523946

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024312849180745
$ws.Range("D16").Value = 1.034606468062034
$ws.Range("E16").Value = 1.024890004573768
$ws.Range("F16").Value = 1.044972499155869
$ws.Range("I16").Value = 1.034981739058093
$ws.Range("J16").Value = 1.030710951335734
$ws.Range("K16").Value = 1.038048954516712
$ws.Range("L16").Value = 1.028367586860331
$ws.Range("M16").Value = 1.048378336606439
$ws.Range("N16").Value = 1.032174679110829

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024683909953545
$ws.Range("D17").Value = 1.034889934210781
$ws.Range("E17").Value = 1.025203873508281
$ws.Range("F17").Value = 1.045317857674197
$ws.Range("I17").Value = 1.035061011927406
$ws.Range("J17").Value = 1.030942488998566
$ws.Range("K17").Value = 1.038259647698994
$ws.Range("L17").Value = 1.028607775402801
$ws.Range("M17").Value = 1.048651543484672
$ws.Range("N17").Value = 1.0324065455837

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024900387590967
$ws.Range("D18").Value = 1.03505527453395
$ws.Range("E18").Value = 1.025387033956001
$ws.Range("F18").Value = 1.045519332021486
$ws.Range("I18").Value = 1.035107055407941
$ws.Range("J18").Value = 1.031077515000649
$ws.Range("K18").Value = 1.038382461016081
$ws.Range("L18").Value = 1.028747881821447
$ws.Range("M18").Value = 1.048810857458184
$ws.Range("N18").Value = 1.032541763338194

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024974208337931
$ws.Range("D19").Value = 1.03511165115323
$ws.Range("E19").Value = 1.025449501484576
$ws.Range("F19").Value = 1.045588035047814
$ws.Range("I19").Value = 1.035122721942675
$ws.Range("J19").Value = 1.031123550940514
$ws.Range("K19").Value = 1.038424323506161
$ws.Range("L19").Value = 1.028795655902417
$ws.Range("M19").Value = 1.048865172078984
$ws.Range("N19").Value = 1.032587864654374

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024644094070074
$ws.Range("D20").Value = 1.034859521036664
$ws.Range("E20").Value = 1.025170189436416
$ws.Range("F20").Value = 1.045280800594777
$ws.Range("I20").Value = 1.035052526868297
$ws.Range("J20").Value = 1.030917649881264
$ws.Range("K20").Value = 1.038237050628534
$ws.Range("L20").Value = 1.028582004558455
$ws.Range("M20").Value = 1.04862223540852
$ws.Range("N20").Value = 1.032381671192003

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023571381585562
$ws.Range("D21").Value = 1.034039816689776
$ws.Range("E21").Value = 1.02426313275234
$ws.Range("F21").Value = 1.044282337655694
$ws.Range("I21").Value = 1.034822024955354
$ws.Range("J21").Value = 1.030247942090262
$ws.Range("K21").Value = 1.037627265145602
$ws.Range("L21").Value = 1.027887505554428
$ws.Range("M21").Value = 1.047831922475748
$ws.Range("N21").Value = 1.031711012339121

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02289732812061
$ws.Range("D22").Value = 1.033524439307741
$ws.Range("E22").Value = 1.023693613950015
$ws.Range("F22").Value = 1.043654867567047
$ws.Range("I22").Value = 1.034675346540566
$ws.Range("J22").Value = 1.02982664221348
$ws.Range("K22").Value = 1.037243151091368
$ws.Range("L22").Value = 1.027450929824206
$ws.Range("M22").Value = 1.047334642687411
$ws.Range("N22").Value = 1.03128911416819

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023254618103715
$ws.Range("D23").Value = 1.033797649841467
$ws.Range("E23").Value = 1.023995451923141
$ws.Range("F23").Value = 1.043987472119828
$ws.Range("I23").Value = 1.034753269466504
$ws.Range("J23").Value = 1.030050002560685
$ws.Range("K23").Value = 1.037446844908666
$ws.Range("L23").Value = 1.027682358651582
$ws.Range("M23").Value = 1.047598295595509
$ws.Range("N23").Value = 1.031512791712708

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024662085017513
$ws.Range("D24").Value = 1.034873263443172
$ws.Range("E24").Value = 1.025185409552949
$ws.Range("F24").Value = 1.045297544993359
$ws.Range("I24").Value = 1.035056361503967
$ws.Range("J24").Value = 1.030928873689584
$ws.Range("K24").Value = 1.038247261520656
$ws.Range("L24").Value = 1.028593649267523
$ws.Range("M24").Value = 1.04863547859947
$ws.Range("N24").Value = 1.032392910939418

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026297101619945
$ws.Range("D25").Value = 1.036121421221948
$ws.Range("E25").Value = 1.026569658363611
$ws.Range("F25").Value = 1.046819083935006
$ws.Range("I25").Value = 1.035400429418843
$ws.Range("J25").Value = 1.031947731490103
$ws.Range("K25").Value = 1.03917294071252
$ws.Range("L25").Value = 1.029651481368691
$ws.Range("M25").Value = 1.049837381331749
$ws.Range("N25").Value = 1.033413215634878
